$d = $word.ActiveDocument

# The document's Title / Author / Abstract paragraphs were each split
# word-by-word across many runs (with separate single-space runs in
# between). Re-set each paragraph's text so Word collapses it back down
# to a single contiguous run, with no change to the visible text.

$targets = @{
    "Title"    = "Answers: Trigonometric identities (degrees)"
    "Author"   = "Dzhemma Ruseva"
    "Abstract" = "A selection of questions on trigonometric identities, using degrees to measure angles."
}

foreach ($styleName in $targets.Keys) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq $styleName) {
            $rng = $d.Range($p.Range.Start, $p.Range.End)
            $rng.Text = $targets[$styleName]
            break
        }
    }
}
